$d = $word.ActiveDocument

# --- 1. Add single-line borders (all sides + inside) to the director-rows table ---
$t = $d.Tables.Item(1)
$t.Borders.Enable = $true

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 2. Rewrite the left cell: split the "line / name" text into a conditional block ---
$leftXml = (
    '<w:p ' + $wNs + '>' +
        '<w:r><w:t xml:space="preserve">{% if </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>row.left</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> %} __________________________  </w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wNs + '>' +
        '<w:r><w:t>{{ row.left.name }} {% endif %}</w:t></w:r>' +
    '</w:p>'
)
$null = $t.Cell(1, 1).Range.InsertXML($leftXml)

# --- 3. Rewrite the right cell: keep the row.right block, then add an else / if-not-row.left block ---
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$rightXml = (
    '<w:p ' + $wNs + '>' +
        '<w:r><w:t xml:space="preserve">{% if </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>row.right</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> %}</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wNs + '>' +
        '<w:r><w:t>__________________________</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wNs + '>' +
        '<w:r><w:t>{{ row.right.name }}</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wNs + '>' +
        '<w:r><w:t>{% else %}</w:t></w:r>' +
    '</w:p>' +
    '<w:p ' + $wNs + '>' +
        '<w:r><w:t xml:space="preserve">{% if not </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>row.left</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> %}{% else %}</w:t></w:r>' +
    '</w:p>'
)
$null = $t.Cell(1, 2).Range.InsertXML($rightXml)

Write-Output "Table updated: borders enabled, left/right cells rewritten."
